$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.0158
$ws.Range("C2").Value = 2.1627
$ws.Range("D2").Value = 6.3036
$ws.Range("E2").Value = 26.248
$ws.Range("F2").Value = 11.0382
$ws.Range("G2").Value = 4.5989

$ws.Range("B3").Value = 1.1352
$ws.Range("C3").Value = 2.0394
$ws.Range("D3").Value = 5.464
$ws.Range("E3").Value = 21.9319
$ws.Range("F3").Value = 8.6128
$ws.Range("G3").Value = 5.7864

$ws.Range("B4").Value = 1.236
$ws.Range("C4").Value = 1.8717
$ws.Range("D4").Value = 4.508
$ws.Range("E4").Value = 17.6459
$ws.Range("F4").Value = 6.746
$ws.Range("G4").Value = 6.26
